$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new weekly data point reported by DGS on 2021/10/27.
# Row 98 (2021/10/25) is the last existing data row; the new data goes
# into row 99.

# Start by cloning the formatting of the previous data row (A98:E98) onto
# the new row so the number formats / styles (date format in col A,
# 0.00 numeric format in cols B:E) carry over exactly.
$ws.Range("A98:E98").Copy()
$ws.Range("A99:E99").PasteSpecial(-4122)

# Write the date label. Prefixing with an apostrophe forces it to be
# stored as literal text (matching how the existing date labels are
# shared-string text, not real date serials) instead of being
# auto-converted into a date value by the date-formatted cell.
$ws.Range("A99").Value = "'2021/10/27"

# Re-apply the source cell's format onto A99 so it keeps style index 1
# (the apostrophe text-entry nudges a cell toward a generic "text"
# style); this restores the yyyy/mm/dd-formatted style used by the rest
# of column A.
$ws.Range("A98").Copy()
$ws.Range("A99").PasteSpecial(-4122)

# Fill in the numeric columns for the new report.
$ws.Range("B99").Value = 94.8
$ws.Range("C99").Value = 94.9
$ws.Range("D99").Value = 1.08
$ws.Range("E99").Value = 1.08

# Move the active selection to the next empty row below the data, same
# as the original sheet's convention of selecting the row right after
# the last populated one.
$ws.Range("A100").Select() | Out-Null
